$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" '34.337.32'
Set-TextValue "E2" '  +12.28%  '
Set-TextValue "D3" '1.827.61'
Set-TextValue "E3" '  +9.47%  '
Set-TextValue "E4" '  -0.13%  '
Set-TextValue "D5" '230.09'
Set-TextValue "E6" '  +8.84%  '
Set-TextValue "E7" '  +0.00%  '
Set-TextValue "D8" '31.55'
Set-TextValue "E8" '  +8.39%  '
Set-TextValue "D9" '46.82'
Set-TextValue "E9" '  +5.98%  '
Set-TextValue "E10" '  +9.90%  '
Set-TextValue "E11" '  +6.27%  '
Set-TextValue "E12" '  +3.28%  '
Set-TextValue "D13" '2.088.15'
Set-TextValue "E13" '  +9.30%  '
Set-TextValue "D14" '1.824.67'
Set-TextValue "E14" '  +9.15%  '
Set-TextValue "E15" '  +8.66%  '
Set-TextValue "D16" '34.312.78'
Set-TextValue "E16" '  +12.16%  '
Set-TextValue "D17" '10.33'
Set-TextValue "E17" '  +3.48%  '
Set-TextValue "E18" '  +7.59%  '
Set-TextValue "D19" '70.58'
Set-TextValue "E19" '  +7.15%  '
Set-TextValue "D20" '258.38'
Set-TextValue "E20" '  +6.79%  '
Set-TextValue "D21" '0.0₃0756'
Set-TextValue "E21" '  +5.35%  '
Set-TextValue "D22" '0.999'
Set-TextValue "E22" '  -0.20%  '
Set-TextValue "D23" '10.66'
Set-TextValue "E23" '  +7.20%  '
Set-TextValue "D24" '4.35'
Set-TextValue "E24" '  +3.04%  '
Set-TextValue "E25" '  +3.62%  '
Set-TextValue "D26" '159.67'
Set-TextValue "E26" '  +0.43%  '
Set-TextValue "D27" '16.77'
Set-TextValue "E27" '  +6.38%  '
Set-TextValue "E28" '  +5.52%  '
Set-TextValue "D29" '7.17'
Set-TextValue "D30" '1.00'
Set-TextValue "E30" '  -0.14%  '
Set-TextValue "E31" '  +13.24%  '
Set-TextValue "E32" '  +6.86%  '
Set-TextValue "E33" '  +6.60%  '
Set-TextValue "D34" '3.58'
Set-TextValue "E34" '  +8.87%  '
Set-TextValue "D35" '1.534.30'
Set-TextValue "E35" '  +2.47%  '
Set-TextValue "E36" '  +3.27%  '
Set-TextValue "E37" '  +6.38%  '
Set-TextValue "D38" '0.637'
Set-TextValue "E38" '  +7.04%  '
Set-TextValue "E39" '  +7.73%  '
Set-TextValue "D40" '84.44'
Set-TextValue "E40" '  +1.27%  '
Set-TextValue "E41" '  +5.44%  '
Set-TextValue "E42" '  +2.94%  '
Set-TextValue "D43" '0.915'
Set-TextValue "E43" '  +9.62%  '
Set-TextValue "D44" '2.13'
Set-TextValue "E44" '  +6.07%  '
Set-TextValue "D45" '0.0527'
Set-TextValue "E45" '  +6.02%  '
Set-TextValue "E46" '  +6.55%  '
Set-TextValue "D47" '1.980.16'
Set-TextValue "E47" '  +9.64%  '
Set-TextValue "E48" '  +5.84%  '
Set-TextValue "D49" '12.17'
Set-TextValue "E49" '  +19.15%  '
Set-TextValue "E50" '  +0.03%  '
Set-TextValue "D51" '51.77'
Set-TextValue "E51" '  +4.23%  '

Write-Host "Applied cryptos list update."
